$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values per regenerated save_data
$kValues = @{
    2 = 1
    3 = 2
    4 = 1
    5 = 0
    6 = 1
    7 = 1
    8 = 0
    9 = 0
    10 = 0
    11 = 2
    12 = 2
    13 = 2
    14 = 0
    15 = 0
    16 = 0
    17 = 2
    18 = 2
    19 = 0
    20 = 2
    21 = 1
    22 = 2
    23 = 1
    24 = 0
    25 = 0
    26 = 1
    27 = 0
    28 = 1
    29 = 2
    30 = 0
    31 = 1
    32 = 0
    33 = 0
    34 = 2
    35 = 2
    36 = 1
    37 = 2
    38 = 1
    39 = 1
    40 = 0
    41 = 1
    42 = 0
    43 = 0
    44 = 0
    45 = 2
    46 = 1
    47 = 1
    48 = 2
    49 = 1
    50 = 2
    51 = 1
    52 = 2
    53 = 1
    54 = 2
    55 = 1
    56 = 1
    57 = 0
    58 = 0
    59 = 1
    60 = 1
    61 = 2
    62 = 1
    65 = 1
    66 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}

